# Daily attendance processing - 2026-01-26 21:08:05
#
# Normalizes the "Recorded By" column (column G) on the active sheet:
# within each comma-separated list of recorders, any entry matching
# "system" (case-insensitively) is moved to the front of the list,
# ahead of the other recorder names/emails, while the relative order
# within each group (system-entries vs. other entries) is preserved.
#
# e.g. "dnasr281@gmail.com, System"              -> "System, dnasr281@gmail.com"
#      "System, backup@backdoor.com, system"     -> "System, system, backup@backdoor.com"
#      "dnasr281@gmail.com, admin@admin.com"     -> unchanged (no "system" entry)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$rowCount = $used.Rows.Count

# Column 7 = G = "Recorded By"
$col = 7

for ($r = 2; $r -le $rowCount; $r++) {
    $cell = $ws.Cells.Item($r, $col)
    $val = $cell.Value()

    if ($val -eq $null -or $val -eq "") {
        continue
    }

    $parts = $val -split ", "

    $systemParts = @()
    $otherParts = @()
    foreach ($p in $parts) {
        if ($p.ToLower() -eq "system") {
            $systemParts += $p
        } else {
            $otherParts += $p
        }
    }

    $newVal = ($systemParts + $otherParts) -join ", "

    if ($newVal -ne $val) {
        $cell.Value = $newVal
    }
}
